# Auto-generated script to update cryptos.xlsx price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.868.32"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.115.79"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "  +0.90%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.34"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = "  -0.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.86"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "  +2.00%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.113.28"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "  +0.92%  "

$ws.Range("E9").Value = "  -0.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.43"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "  -2.87%  "

$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("E13").Value = "  -1.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.97"
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = "  +1.21%  "

$ws.Range("E15").Value = "  -1.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.632.70"
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = "  +0.96%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.830.62"
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = "  -0.15%  "

$ws.Range("E18").Value = "  -0.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.113.14"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "  +0.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.31"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = "  +1.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "477.34"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "  +2.50%  "

$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.04"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = "  +0.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.31"
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = "  +3.18%  "

$ws.Range("E26").Value = "  -3.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.08"
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = "  -0.69%  "

$ws.Range("E30").Value = "  -1.30%  "

$ws.Range("E31").Value = "  +0.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.56"
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = "  +1.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.116"
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = "  +1.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0942"
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = "  -7.44%  "

$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("E36").Value = "  -0.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.974"
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = "  -3.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.23"
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = "  +1.05%  "

$ws.Range("E39").Value = "  -0.06%  "

$ws.Range("E40").Value = "  -3.63%  "

$ws.Range("E41").Value = "  -2.42%  "

$ws.Range("E42").Value = "  +0.35%  "

$ws.Range("E43").Value = "  -0.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.825.16"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "383.34"
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = "  -0.15%  "

$ws.Range("E46").Value = "  -1.67%  "

$ws.Range("E47").Value = "  -9.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.29"
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = "  +0.50%  "

$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("E50").Value = "  +0.56%  "

$ws.Range("E51").Value = "  -1.97%  "
